$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.970.16"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "1.653.10"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'214.96"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.248"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "'0.0613"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'19.50"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "1.891.20"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("D13").Value = "1.675.69"
$ws.Range("E13").Value = "  +7.11%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "'64.70"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "27.035.41"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "'237.04"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "'7.81"
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("D20").Value = "0.0₃0727"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'4.42"
$ws.Range("E22").Value = "  +4.51%  "
$ws.Range("D23").Value = "'2.24"
$ws.Range("E23").Value = "  +5.16%  "
$ws.Range("D24").Value = "'9.25"
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("D25").Value = "'146.06"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'7.14"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "'15.76"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").Value = "'0.0496"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "1.531.61"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'0.572"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0168"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.881"
$ws.Range("E39").Value = "  +7.58%  "
$ws.Range("D40").Value = "'5.94"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'2.27"
$ws.Range("E42").Value = "  +4.38%  "
$ws.Range("D43").Value = "'66.12"
$ws.Range("E43").Value = "  +9.46%  "
$ws.Range("D44").Value = "1.797.93"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("D45").Value = "'0.775"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").Value = "'0.922"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'89.99"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "'1.52"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'0.0973"
$ws.Range("E51").Value = "  +3.01%  "
